# Revert column widths and selection on the "Combined" sheet
# (undo the merge of branch 'scraper').
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined")

# Restore original column widths for columns A and D.
# Note: the runtime's ColumnWidth setter stores width + 0.8333333333333334
# (5/6 of a character) in the XML "width" attribute, so subtract that
# offset here to land exactly on the target stored widths of 17 and 20.
$ws.Columns.Item(1).ColumnWidth = 17 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 20 - 0.8333333333333334

# Move the active cell selection to C11.
$ws.Activate()
$ws.Range("C11").Select()
